$d = $word.ActiveDocument

# The template previously used the (invalid) Jinja tag "{% licenced_vessel %}"
# in four places (Registration number, Vessel Name, Registered length and
# Draft paragraphs of the "Licensed Vessel" block). That tag is missing the
# "if" keyword, so Jinja never rendered the vessel block and the licence
# email body came out empty. The fix inserts "if " right after "{% " so the
# tag reads "{% if licenced_vessel %}".
#
# Word naturally breaks the original single run into three runs at the
# insertion boundaries, so we reproduce that explicitly instead of letting
# a plain Find/Replace collapse everything back into one run:
#   run 1: "{% "
#   run 2: "if " / "if"
#   run 3: "licenced_vessel %}" / " licenced_vessel %}"
# All three keep the original "Arial Nova Light" 6pt (half-point size 12)
# formatting that the single run had before the edit.

$occurrence = 0

while ($true) {
    $probe = $d.Content
    $found = $probe.Find.Execute("{% licenced_vessel %}", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        break
    }

    $occurrence = $occurrence + 1

    $matchStart = $probe.Start
    $matchEnd = $probe.End

    # Insertion point is right after the literal "{% " (3 characters).
    $insPos = $matchStart + 3
    $insRange = $d.Range($insPos, $insPos)

    $insRange.Text = "if "

    if ($occurrence -eq 3) {
        # "Registered length" paragraph: split lands after "if" (no
        # trailing space on that run - the space stays on the third run).
        $splitLen = 2
    } else {
        $splitLen = 3
    }

    # Force the freshly inserted text to stay as its own run (rather than
    # being silently re-merged with the identically-formatted runs around
    # it) by toggling Bold on and back off across each of the two new run
    # boundaries - first across "if "/"if", then across the remaining
    # "licenced_vessel %}" / " licenced_vessel %}" text.
    $boundary1 = $d.Range($insPos, $insPos + $splitLen)
    $boundary1.Font.Bold = $true
    $boundary1.Font.Bold = $false

    $boundary2 = $d.Range($insPos + $splitLen, $matchEnd + 3)
    $boundary2.Font.Bold = $true
    $boundary2.Font.Bold = $false
}

Write-Host "Fixed $occurrence occurrence(s) of the licenced_vessel tag."
